{"js": "// The \"Requisitos\" (requirements) list at the end of the document is a\n// single paragraph made up of one run per course (\"<course> - <name>\n// (Requisito)\") each followed by a manual line break. The edit re-orders\n// those course lines (same 54 lines, new order) without touching anything\n// else in the document.\n\nconst FINAL_ORDER = [\n  \"LOQ4204 -  Economia Geral  (Requisito)\",\n  \"LOQ4240 -  Administra\u00e7\u00e3o e Organiza\u00e7\u00e3o II  (Requisito)\",\n  \"LOB1019 -  F\u00edsica II  (Requisito)\",\n  \"LOQ4203 -  Sistemas Produtivos  (Requisito)\",\n  \"LOB1006 -  C\u00e1lculo IV  (Requisito)\",\n  \"LOQ4209 -  Engenharia da Qualidade  (Requisito)\",\n  \"LOQ4095 -  Qu\u00edmica Geral Experimental  (Requisito)\",\n  \"LOQ4201 -  Introdu\u00e7\u00e3o \u00e0 Engenharia de Produ\u00e7\u00e3o  (Requisito)\",\n  \"LOQ4236 -  Projeto Integrado de Engenharia de Produ\u00e7\u00e3o I  (Requisito)\",\n  \"LOQ4239 -  Administra\u00e7\u00e3o e Organiza\u00e7\u00e3o I  (Requisito)\",\n  \"LOQ4251 -  Fundamentos de Qu\u00edmica  (Requisito)\",\n  \"LOB1004 -  C\u00e1lculo II  (Requisito)\",\n  \"LOB1018 -  F\u00edsica I  (Requisito)\",\n  \"LOQ4257 -  Gest\u00e3o de Projetos  (Requisito)\",\n  \"LOB1039 -  F\u00edsica Experimental III  (Requisito)\",\n  \"LOQ4264 -  Engenharia da Sustentabilidade  (Requisito)\",\n  \"LOQ4262 -  Automa\u00e7\u00e3o e Controle  (Requisito)\",\n  \"LOQ4271 -  Planejamento de Experimentos  (Requisito)\",\n  \"LOB1003 -  C\u00e1lculo I  (Requisito)\",\n  \"LOB1009 -  Leitura e Interpreta\u00e7\u00e3o de Desenho T\u00e9cnico  (Requisito)\",\n  \"LOB1036 -  Geometria Anal\u00edtica  (Requisito)\",\n  \"LOB1038 -  F\u00edsica Experimental I  (Requisito)\",\n  \"LOM3016 -  Introdu\u00e7\u00e3o \u00e0  Ci\u00eancia dos Materiais  (Requisito)\",\n  \"LOB1012 -  Estat\u00edstica  (Requisito)\",\n  \"LOB1041 -  F\u00edsica Experimental II  (Requisito)\",\n  \"LOB1052 -  C\u00e1lculo III  (Requisito)\",\n  \"LOB1049 -  Estat\u00edstica Multivariada  (Requisito)\",\n  \"LOB1053 -  F\u00edsica III  (Requisito)\",\n  \"LOQ4237 -  Projeto Integrado de Engenharia de Produ\u00e7\u00e3o II  (Requisito)\",\n  \"LOB1024 -  Mec\u00e2nica  (Requisito)\",\n  \"LOB1046 -  Engenharia do Meio Ambiente  (Requisito)\",\n  \"LOB1056 -  Introdu\u00e7\u00e3o aos M\u00e9todos Num\u00e9ricos e Computacionais  (Requisito)\",\n  \"LOQ4206 -  Pesquisa Operacional I  (Requisito)\",\n  \"LOQ4252 -  Fundamentos de Fen\u00f4menos de Transporte  (Requisito)\",\n  \"LOQ4253 -  Processos Qu\u00edmicos Industriais  (Requisito)\",\n  \"LOB1011 -  Eletricidade Aplicada  (Requisito)\",\n  \"LOB1040 -  Laborat\u00f3rio de Eletricidade  (Requisito)\",\n  \"LOM3081 -  Introdu\u00e7\u00e3o \u00e0 Mec\u00e2nica dos S\u00f3lidos  (Requisito)\",\n  \"LOQ4076 -  Termodin\u00e2mica Aplicada  (Requisito)\",\n  \"LOQ4234 -  Empreendedorismo e Inova\u00e7\u00e3o  (Requisito)\",\n  \"LOQ4255 -  Inova\u00e7\u00e3o Tecnol\u00f3gica  (Requisito)\",\n  \"LOQ4258 -  Pesquisa Operacional II  (Requisito)\",\n  \"LOQ4263 -  Planejamento e Gest\u00e3o da Manuten\u00e7\u00e3o  (Requisito)\",\n  \"LOQ4213 -  Contabilidade e Custos  (Requisito)\",\n  \"LOQ4238 -  Projeto Integrado em Engenharia de Produ\u00e7\u00e3o III  (Requisito)\",\n  \"LOQ4241 -  Sistemas de Apoio \u00e0 Decis\u00e3o  (Requisito)\",\n  \"LOQ4245 -  Ergonomia  (Requisito)\",\n  \"LOQ4259 -  Processos de Desenvolvimento de Servi\u00e7os  (Requisito)\",\n  \"LOQ4261 -  Planejamento, Programa\u00e7\u00e3o e Controle da Produ\u00e7\u00e3o I  (Requisito)\",\n  \"LOB1055 -  Fundamentos de Engenharia de Seguran\u00e7a no Trabalho  (Requisito)\",\n  \"LOQ4222 -  Engenharia Econ\u00f4mica e Finan\u00e7as  (Requisito)\",\n  \"LOQ4260 -  Controle Estat\u00edstico da Qualidade  (Requisito)\",\n  \"LOQ4270 -  Planejamento, Programa\u00e7\u00e3o e Controle da Produ\u00e7\u00e3o II  (Requisito)\",\n  \"LOQ4272 -  Projeto da Fabrica  (Requisito)\"\n];\n\nfunction xmlEscape(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\")\n    .replace(/'/g, \"&apos;\");\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"Requisitos\" heading paragraph, the requirements list is the\n// very next paragraph (style \"ListBullet\").\nlet headingIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Requisitos\") {\n    headingIndex = i;\n    break;\n  }\n}\n\nif (headingIndex === -1 || headingIndex + 1 >= paragraphs.items.length) {\n  throw new Error(\"Could not locate the 'Requisitos' list paragraph.\");\n}\n\nconst listParagraph = paragraphs.items[headingIndex + 1];\nlistParagraph.load(\"style\");\nawait context.sync();\n\n// Rebuild the runs of the requirements paragraph in the new order, each\n// course on its own run followed by a manual line break, exactly like the\n// original markup (<w:r><w:t>\u2026</w:t><w:br/></w:r> repeated).\nconst runsXml = FINAL_ORDER.map(\n  (line) => `<w:r><w:t>${xmlEscape(line)}</w:t><w:br/></w:r>`\n).join(\"\");\n\nconst ooxml =\n  `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>` +\n  `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">` +\n  `<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>` +\n  `<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">` +\n  `<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>` +\n  `</Relationships>` +\n  `</pkg:xmlData></pkg:part>` +\n  `<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>` +\n  `<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">` +\n  `<w:body><w:p><w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr>${runsXml}</w:p></w:body>` +\n  `</w:document>` +\n  `</pkg:xmlData></pkg:part>` +\n  `</pkg:package>`;\n\nlistParagraph.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The \"Requisitos\" (requirements) list at the end of the document is a\n# single paragraph made up of one run per course (\"<course> - <name>\n# (Requisito)\") each followed by a manual line break. The edit re-orders\n# those course lines (same 54 lines, new order) without touching anything\n# else in the document.\n\n$FinalOrder = @(\n    \"LOQ4204 -  Economia Geral  (Requisito)\",\n    \"LOQ4240 -  Administra\u00e7\u00e3o e Organiza\u00e7\u00e3o II  (Requisito)\",\n    \"LOB1019 -  F\u00edsica II  (Requisito)\",\n    \"LOQ4203 -  Sistemas Produtivos  (Requisito)\",\n    \"LOB1006 -  C\u00e1lculo IV  (Requisito)\",\n    \"LOQ4209 -  Engenharia da Qualidade  (Requisito)\",\n    \"LOQ4095 -  Qu\u00edmica Geral Experimental  (Requisito)\",\n    \"LOQ4201 -  Introdu\u00e7\u00e3o \u00e0 Engenharia de Produ\u00e7\u00e3o  (Requisito)\",\n    \"LOQ4236 -  Projeto Integrado de Engenharia de Produ\u00e7\u00e3o I  (Requisito)\",\n    \"LOQ4239 -  Administra\u00e7\u00e3o e Organiza\u00e7\u00e3o I  (Requisito)\",\n    \"LOQ4251 -  Fundamentos de Qu\u00edmica  (Requisito)\",\n    \"LOB1004 -  C\u00e1lculo II  (Requisito)\",\n    \"LOB1018 -  F\u00edsica I  (Requisito)\",\n    \"LOQ4257 -  Gest\u00e3o de Projetos  (Requisito)\",\n    \"LOB1039 -  F\u00edsica Experimental III  (Requisito)\",\n    \"LOQ4264 -  Engenharia da Sustentabilidade  (Requisito)\",\n    \"LOQ4262 -  Automa\u00e7\u00e3o e Controle  (Requisito)\",\n    \"LOQ4271 -  Planejamento de Experimentos  (Requisito)\",\n    \"LOB1003 -  C\u00e1lculo I  (Requisito)\",\n    \"LOB1009 -  Leitura e Interpreta\u00e7\u00e3o de Desenho T\u00e9cnico  (Requisito)\",\n    \"LOB1036 -  Geometria Anal\u00edtica  (Requisito)\",\n    \"LOB1038 -  F\u00edsica Experimental I  (Requisito)\",\n    \"LOM3016 -  Introdu\u00e7\u00e3o \u00e0  Ci\u00eancia dos Materiais  (Requisito)\",\n    \"LOB1012 -  Estat\u00edstica  (Requisito)\",\n    \"LOB1041 -  F\u00edsica Experimental II  (Requisito)\",\n    \"LOB1052 -  C\u00e1lculo III  (Requisito)\",\n    \"LOB1049 -  Estat\u00edstica Multivariada  (Requisito)\",\n    \"LOB1053 -  F\u00edsica III  (Requisito)\",\n    \"LOQ4237 -  Projeto Integrado de Engenharia de Produ\u00e7\u00e3o II  (Requisito)\",\n    \"LOB1024 -  Mec\u00e2nica  (Requisito)\",\n    \"LOB1046 -  Engenharia do Meio Ambiente  (Requisito)\",\n    \"LOB1056 -  Introdu\u00e7\u00e3o aos M\u00e9todos Num\u00e9ricos e Computacionais  (Requisito)\",\n    \"LOQ4206 -  Pesquisa Operacional I  (Requisito)\",\n    \"LOQ4252 -  Fundamentos de Fen\u00f4menos de Transporte  (Requisito)\",\n    \"LOQ4253 -  Processos Qu\u00edmicos Industriais  (Requisito)\",\n    \"LOB1011 -  Eletricidade Aplicada  (Requisito)\",\n    \"LOB1040 -  Laborat\u00f3rio de Eletricidade  (Requisito)\",\n    \"LOM3081 -  Introdu\u00e7\u00e3o \u00e0 Mec\u00e2nica dos S\u00f3lidos  (Requisito)\",\n    \"LOQ4076 -  Termodin\u00e2mica Aplicada  (Requisito)\",\n    \"LOQ4234 -  Empreendedorismo e Inova\u00e7\u00e3o  (Requisito)\",\n    \"LOQ4255 -  Inova\u00e7\u00e3o Tecnol\u00f3gica  (Requisito)\",\n    \"LOQ4258 -  Pesquisa Operacional II  (Requisito)\",\n    \"LOQ4263 -  Planejamento e Gest\u00e3o da Manuten\u00e7\u00e3o  (Requisito)\",\n    \"LOQ4213 -  Contabilidade e Custos  (Requisito)\",\n    \"LOQ4238 -  Projeto Integrado em Engenharia de Produ\u00e7\u00e3o III  (Requisito)\",\n    \"LOQ4241 -  Sistemas de Apoio \u00e0 Decis\u00e3o  (Requisito)\",\n    \"LOQ4245 -  Ergonomia  (Requisito)\",\n    \"LOQ4259 -  Processos de Desenvolvimento de Servi\u00e7os  (Requisito)\",\n    \"LOQ4261 -  Planejamento, Programa\u00e7\u00e3o e Controle da Produ\u00e7\u00e3o I  (Requisito)\",\n    \"LOB1055 -  Fundamentos de Engenharia de Seguran\u00e7a no Trabalho  (Requisito)\",\n    \"LOQ4222 -  Engenharia Econ\u00f4mica e Finan\u00e7as  (Requisito)\",\n    \"LOQ4260 -  Controle Estat\u00edstico da Qualidade  (Requisito)\",\n    \"LOQ4270 -  Planejamento, Programa\u00e7\u00e3o e Controle da Produ\u00e7\u00e3o II  (Requisito)\",\n    \"LOQ4272 -  Projeto da Fabrica  (Requisito)\"\n)\n\n$d = $word.ActiveDocument\n\n# Locate the \"Requisitos\" heading paragraph; the requirements list is the\n# very next paragraph (style \"ListBullet\").\n$listParagraph = $null\nforeach ($para in $d.Paragraphs) {\n    $t = $para.Range.Text.Trim()\n    if ($t -eq \"Requisitos\") {\n        $listParagraph = $para.Next()\n        break\n    }\n}\n\nif ($listParagraph -eq $null) {\n    throw \"Could not locate the 'Requisitos' list paragraph.\"\n}\n\n$startPos = $listParagraph.Range.Start\n$endPos = $listParagraph.Range.End\n$targetRange = $d.Range($startPos, $endPos)\n\n# Rebuild the runs of the requirements paragraph in the new order, each\n# course on its own run followed by a manual line break, exactly like the\n# original markup (<w:r><w:t>\u2026</w:t><w:br/></w:r> repeated).\nfunction XmlEscape($s) {\n    $s = $s -replace \"&\", \"&amp;\"\n    $s = $s -replace \"<\", \"&lt;\"\n    $s = $s -replace \">\", \"&gt;\"\n    $s = $s -replace '\"', \"&quot;\"\n    $s = $s -replace \"'\", \"&apos;\"\n    return $s\n}\n\n$runsXml = \"\"\nforeach ($line in $FinalOrder) {\n    $escaped = XmlEscape $line\n    $runsXml += '<w:r><w:t>' + $escaped + '</w:t><w:br/></w:r>'\n}\n\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' `\n    + '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' `\n    + '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' `\n    + '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' `\n    + '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' `\n    + '</Relationships>' `\n    + '</pkg:xmlData></pkg:part>' `\n    + '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' `\n    + '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' `\n    + '<w:body><w:p><w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr>' + $runsXml + '</w:p></w:body>' `\n    + '</w:document>' `\n    + '</pkg:xmlData></pkg:part>' `\n    + '</pkg:package>'\n\n$targetRange.InsertXML($ooxml)\n"}
